$d = $word.ActiveDocument

# Remove trailing ",00" / ",0" decimals and update amounts per the
# "suppression des 0,00" cleanup commit.

$d.Content.Find.Execute("8140,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "8140", 2)

$d.Content.Find.Execute("25505,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "25505", 2)

$d.Content.Find.Execute("6405,0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "6100", 2)

$d.Content.Find.Execute("14565,0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "14260", 2)

$d.Content.Find.Execute("57110,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "56805", 2)
